# Move the automatic "_GoBack" bookmark (Word's "last edit location"
# marker) from where it currently sits -- right at the start of the
# Title paragraph -- to the point inside the Content paragraph where the
# user's last edit happened, i.e. right after "[Con" in "[Content]".
# This also has the effect of splitting the single run that holds
# "[Content]" into two runs: "[Con" and "tent]".
$d = $word.ActiveDocument

$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$contentBookmark = $d.Bookmarks("Content")
$splitPoint = $contentBookmark.Range.Start + 4   # after "[Con"
$r = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $r)

# Add a new inline-code character style, based on the existing
# CodeSnippet character style, for marking up inline code spans.
$inlineCode = $d.Styles.Add("InlineCodeSnippet", 2)
$inlineCode.BaseStyle = $d.Styles("CodeSnippetZchn")
$inlineCode.Priority = 1
$inlineCode.QuickStyle = $true
